# Modifica os artefatos de front
# Inserts a new "Properties" worksheet between "Metadata" and "Concepts",
# describing the "category" property (Code/Uri/Description/Type columns).

$wb = $excel.ActiveWorkbook

$metadataSheet = $wb.Worksheets.Item("Metadata")

# Insert the new sheet right after "Metadata" (i.e. before the existing "Concepts" sheet)
$propsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $metadataSheet)
$propsSheet.Name = "Properties"

# Header row content
$propsSheet.Range("A1").Value = "Code"
$propsSheet.Range("B1").Value = "Uri"
$propsSheet.Range("C1").Value = "Description"
$propsSheet.Range("D1").Value = "Type"

# Data row content (the "category" property)
$propsSheet.Range("A2").Value = "category"
$propsSheet.Range("B2").Value = "http://www.saude.gov.br/fhir/r4/ValueSet/BRCategoriaAgenteAlergiasReacoesAdversas-1.0"
$propsSheet.Range("D2").Value = "code"

# Match the look of the existing sheets: bold header style copied from the
# "Metadata" sheet's header row, plain style copied from its body rows.
$metadataSheet.Range("A1:B1").Copy()
$propsSheet.Range("A1:D1").PasteSpecial(-4122)

$metadataSheet.Range("A2:B2").Copy()
$propsSheet.Range("A2:D2").PasteSpecial(-4122)
